$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Occupied Bed Numbers" column (column J) first, then the
# "Daily Rent" column (column G), so column indices don't shift unexpectedly.
$ws.Columns.Item(10).Delete()
$ws.Columns.Item(7).Delete()
